# TemplateNhapDiemDotThi.xlsx - "fix nhap diem :(" edit
#
# 1) Data fix: the 4 score cells (msword/msexcel/mspowerpoint + sobaodanh date
#    cell) were being stored as numbers; re-enter them as text "0" / "36483"
#    so they match the "text" number format (style 1, numFmtId 49) already
#    applied to these columns.
# 2) View fix: scroll/select over to column B and zoom to 85%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table layout: set the "0" cells first,
# then the date-like "36483" cell, matching how the sheet was actually edited.
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "0"
$ws.Range("L2").Value = "0"
$ws.Range("C2").Value = "36483"

# View state: select B2 and zoom to 85%.
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 85
